# Scheduled runner update: refresh market-board derived price/profit
# columns (H:N) across the Leve profit sheets. Only the computed columns
# change; item/leve metadata (A:G) is left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 73.73333
$ws.Range("I11").Value = 73.73333
$ws.Range("K11").Value = 73.73333
$ws.Range("M11").Value = 66.26667

$ws.Range("H38").Value = 1413.0555
$ws.Range("I38").Value = 93.5
$ws.Range("J38").Value = 3062.5
$ws.Range("K38").Value = 280.5
$ws.Range("L38").Value = 9187.5
$ws.Range("M38").Value = 91.5
$ws.Range("N38").Value = -9931.5

$ws.Range("H70").Value = 2007.0714
$ws.Range("I70").Value = 2960
$ws.Range("J70").Value = 1477.6666
$ws.Range("K70").Value = 8880
$ws.Range("L70").Value = 4432.9998
$ws.Range("M70").Value = -8610
$ws.Range("N70").Value = -4972.9998

$ws.Range("H73").Value = 2007.0714
$ws.Range("I73").Value = 2960
$ws.Range("J73").Value = 1477.6666
$ws.Range("K73").Value = 8880
$ws.Range("L73").Value = 4432.9998
$ws.Range("M73").Value = -7944
$ws.Range("N73").Value = -6304.9998

$ws.Range("H125").Value = 789.375
$ws.Range("I125").Value = 789.375
$ws.Range("K125").Value = 7104.375
$ws.Range("M125").Value = -4644.375

$ws.Range("H129").Value = 1048.0646
$ws.Range("I129").Value = 284.83334
$ws.Range("J129").Value = 1231.24
$ws.Range("K129").Value = 854.5000200000001
$ws.Range("L129").Value = 3693.72
$ws.Range("M129").Value = 4145.49998
$ws.Range("N129").Value = -13693.72

$ws.Range("H131").Value = 535.2222
$ws.Range("I131").Value = 227.125
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 681.375
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = 4358.625
$ws.Range("N131").Value = -19080

$ws.Range("H132").Value = 8151.25
$ws.Range("I132").Value = 8151.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24453.75
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -21923.75

$ws.Range("H137").Value = 1143.3191
$ws.Range("I137").Value = 846.7222
$ws.Range("J137").Value = 2114
$ws.Range("K137").Value = 2540.1666
$ws.Range("L137").Value = 6342
$ws.Range("M137").Value = 9.833399999999983
$ws.Range("N137").Value = -11442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27778692
$ws.Range("I2").Value = 35715116
$ws.Range("J2").Value = 1209.25
$ws.Range("K2").Value = 35715116
$ws.Range("L2").Value = 1209.25
$ws.Range("M2").Value = -35715003
$ws.Range("N2").Value = -1435.25

$ws.Range("H88").Value = 2698.6667
$ws.Range("I88").Value = 1837.2
$ws.Range("J88").Value = 3314
$ws.Range("K88").Value = 1837.2
$ws.Range("L88").Value = 3314
$ws.Range("M88").Value = -1431.2
$ws.Range("N88").Value = -4126

$ws.Range("H91").Value = 2698.6667
$ws.Range("I91").Value = 1837.2
$ws.Range("J91").Value = 3314
$ws.Range("K91").Value = 1837.2
$ws.Range("L91").Value = 3314
$ws.Range("M91").Value = -433.2
$ws.Range("N91").Value = -6122

$ws.Range("H116").Value = 27778692
$ws.Range("I116").Value = 35715116
$ws.Range("J116").Value = 1209.25
$ws.Range("K116").Value = 35715116
$ws.Range("L116").Value = 1209.25
$ws.Range("M116").Value = -35712822
$ws.Range("N116").Value = -5797.25

$ws.Range("H122").Value = 1503.9445
$ws.Range("I122").Value = 1284.7742
$ws.Range("J122").Value = 2862.8
$ws.Range("K122").Value = 3854.3226
$ws.Range("L122").Value = 8588.400000000001
$ws.Range("M122").Value = -1404.3226
$ws.Range("N122").Value = -13488.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27778692
$ws.Range("I3").Value = 35715116
$ws.Range("J3").Value = 1209.25
$ws.Range("K3").Value = 35715116
$ws.Range("L3").Value = 1209.25
$ws.Range("M3").Value = -35715002
$ws.Range("N3").Value = -1437.25

$ws.Range("H86").Value = 2068.8462
$ws.Range("I86").Value = 1605.625
$ws.Range("K86").Value = 1605.625
$ws.Range("M86").Value = -482.625

$ws.Range("H89").Value = 2068.8462
$ws.Range("I89").Value = 1605.625
$ws.Range("K89").Value = 7325
$ws.Range("M89").Value = -2412.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5007.4614
$ws.Range("I99").Value = 4554.273
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 4554.273
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -3056.273
$ws.Range("N99").Value = -10496

$ws.Range("H126").Value = 5007.4614
$ws.Range("I126").Value = 4554.273
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 13662.819
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -11192.819
$ws.Range("N126").Value = -27440

$ws.Range("H135").Value = 5357934
$ws.Range("J135").Value = 5357934
$ws.Range("L135").Value = 5357934
$ws.Range("N135").Value = -5368074

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3998.5557
$ws.Range("I137").Value = 4010
$ws.Range("J137").Value = 3997.5151
$ws.Range("K137").Value = 12030
$ws.Range("L137").Value = 11992.5453
$ws.Range("M137").Value = -6930
$ws.Range("N137").Value = -22192.5453

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2793.4814
$ws.Range("I80").Value = 2656.8
$ws.Range("J80").Value = 2873.8823
$ws.Range("K80").Value = 2656.8
$ws.Range("L80").Value = 2873.8823
$ws.Range("M80").Value = -1658.8
$ws.Range("N80").Value = -4869.8823

$ws.Range("H83").Value = 2793.4814
$ws.Range("I83").Value = 2656.8
$ws.Range("J83").Value = 2873.8823
$ws.Range("K83").Value = 13284
$ws.Range("L83").Value = 14369.4115
$ws.Range("M83").Value = -8292
$ws.Range("N83").Value = -24353.4115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 42542
$ws.Range("J108").Value = 42542
$ws.Range("L108").Value = 42542
$ws.Range("N108").Value = -50222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2323.52
$ws.Range("I81").Value = 2215.6667
$ws.Range("J81").Value = 2423.077
$ws.Range("K81").Value = 4431.3334
$ws.Range("L81").Value = 4846.154
$ws.Range("M81").Value = -3370.3334
$ws.Range("N81").Value = -6968.154

$ws.Range("H84").Value = 2323.52
$ws.Range("I84").Value = 2215.6667
$ws.Range("J84").Value = 2423.077
$ws.Range("K84").Value = 22156.667
$ws.Range("L84").Value = 24230.77
$ws.Range("M84").Value = -16852.667
$ws.Range("N84").Value = -34838.77
